$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style/format of H1 (bold, border, centered) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows: row, I value, J value
$rows = @(
    @(2, 6, 7),
    @(3, 1, 6),
    @(4, 1, 4),
    @(5, 1, 6),
    @(6, 1, 6),
    @(7, 1, 7),
    @(8, 1, 6),
    @(9, 1, 6),
    @(10, 1, 5),
    @(11, 1, 3),
    @(12, 1, 6),
    @(13, 1, 6),
    @(14, 1, 7),
    @(15, 1, 6),
    @(16, 1, 6),
    @(17, 1, 7),
    @(18, 1, 6),
    @(19, 1, 7),
    @(20, 1, 5),
    @(21, 1, 7),
    @(22, 1, 5),
    @(23, 1, 7),
    @(24, 1, 7),
    @(25, 1, 8),
    @(26, 1, 8),
    @(27, 1, 7),
    @(28, 1, 6),
    @(29, 1, 6),
    @(30, 1, 6),
    @(31, 1, 6),
    @(32, 1, 7),
    @(33, 1, 4),
    @(34, 1, 5),
    @(35, 1, 6),
    @(36, 1, 3),
    @(37, 1, 2)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 9).Value = $r[1]
    $ws.Cells.Item($rowNum, 10).Value = $r[2]
}
